$wb = $excel.ActiveWorkbook

# ---- About sheet (unchanged data; keep as-is) ----
$wsAbout = $wb.Worksheets.Item("About")

# ---- FPIEBP sheet: update production/imports/exports priority values ----
$wsF = $wb.Worksheets.Item("FPIEBP")

# Rows whose Production/Imports/Exports priorities become 2 / 1 / 3
$rowsToSwap = @(3,4,9,10,11,12,13,14,17,18,19,20,22)

foreach ($r in $rowsToSwap) {
    $rng = $wsF.Range("B" + $r + ":D" + $r)
    $rng.ClearFormats()
    $wsF.Range("B" + $r).Value = 2
    $wsF.Range("C" + $r).Value = 1
    $wsF.Range("D" + $r).Value = 3
}

# Row 22's Production cell (B22) keeps its distinct highlight style
$wsF.Range("B22").Style = "Normal"
$wsF.Range("B22").Interior.ColorIndex = -4142

# Page setup for FPIEBP
$wsF.PageSetup.PaperSize = [Microsoft.Office.Interop.Excel.XlPaperSize]::xlPaperA4
$wsF.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait

# Make FPIEBP the active/selected sheet and select D23, matching the saved view state
$wsF.Activate()
$wsF.Range("D23").Select()
